$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 63, shifting rows 63:105 down to 64:106.
$ws.Rows.Item(63).Insert()

# Populate the newly inserted row 63 with the new weekly record.
$ws.Range("A63").Value = 10
$ws.Range("B63").Value = 'Vega Modelo de Temuco'
$ws.Range("C63").Value = 'La Araucanía'
$ws.Range("D63").Value = 45126
$ws.Range("E63").Value = 9
$ws.Range("F63").Value = 'Fruta'
$ws.Range("G63").Value = 100108
$ws.Range("H63").Value = 'Tropicales y subtropicales'
$ws.Range("I63").Value = 100108007
$ws.Range("J63").Value = 'Coco'
$ws.Range("K63").Value = 'Sin especificar'
$ws.Range("L63").Value = 'Primera'
$ws.Range("M63").Value = 20
$ws.Range("N63").Value = 36000
$ws.Range("O63").Value = 36000
$ws.Range("P63").Value = 36000
$ws.Range("Q63").Value = '$/malla 20 unidades'
$ws.Range("R63").Value = 'Perú'
$ws.Range("S63").Value = 1800
$ws.Range("T63").Value = 20
